$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the daily time-series with 4 more days (rows 230-233),
# mirroring the existing columns: A=date serial, B=nuovi pos.,
# C=somma mobile 7gg., D=somma mobile 7gg. per 100mila abitanti.
$data = @(
    @(44304, 2, 3, 140.1214385801028),
    @(44305, 0, 3, 140.1214385801028),
    @(44306, 0, 3, 140.1214385801028),
    @(44307, 0, 3, 140.1214385801028)
)

$row = 230
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row = $row + 1
}

# Reuse the date column's existing formatting (bold, centered, boxed,
# custom date/time number format) by copying it from the row above,
# same as dragging the fill handle down in Excel.
$ws.Range("A229").Copy()
$ws.Range("A230:A233").PasteSpecial(-4122)
$excel.CutCopyMode = $false
